$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 54
$ws.Cells.Item($row, 1).Value = "2025-04-29 06:53:46"
$ws.Cells.Item($row, 2).Value = 146
